$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the numeric annotation scores for columns E:J (Clear, Assertive,
# Cautious, Optimistic, Specific, Relevant) for rows 2-10.
$data = @(
    @(2,1,1,1,2,2),
    @(2,2,1,1,1,2),
    @(2,1,1,1,2,2),
    @(2,0,1,1,2,2),
    @(2,2,1,2,2,2),
    @(2,2,1,1,2,2),
    @(2,1,1,1,2,2),
    @(2,2,1,2,2,2),
    @(2,2,1,2,1,2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 5 + $j  # Column E = 5
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$j]
        # Match the center/wrap formatting already used across the row.
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4108
        $cell.WrapText = $true
    }
}

# Update sheet view: zoom to 85%, freeze the header row, and set the active
# selection to G11 within the frozen (bottom-left) pane.
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("G11").Select()
